$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Seven new province/district rows were added to the source data, which is
# kept sorted alphabetically by Province. Insert blank rows at the positions
# they occupy in the FINAL layout, top-to-bottom: since every row currently
# above the insertion point is already in its final resting place, inserting
# at that row index pushes the remaining (not-yet-finalized) rows down by one
# and reproduces the target arrangement.
$insertPositions = @(27,28,30,55,56,57,58)
foreach ($pos in $insertPositions) {
    $ws.Rows.Item($pos).Insert()
}

# Full target content for rows 27-88: Province, Region, T_Evacuation, T_Health, T_School
$data = @(
    @("CITY OF ISABELA (NOT A PROVINCE)","REGION IX (ZAMBOANGA PENINSULA)",0,0,62),
    @("COMPOSTELA VALLEY","REGION XI (DAVAO REGION)",185,0,404),
    @("COTABATO (NORTH COTABATO)","REGION XII (SOCCSKSARGEN)",305,587,910),
    @("COTABATO CITY (NOT A PROVINCE)","REGION XII (SOCCSKSARGEN)",4,0,43),
    @("DAVAO DEL NORTE","REGION XI (DAVAO REGION)",149,315,425),
    @("DAVAO DEL SUR","REGION XI (DAVAO REGION)",249,595,664),
    @("DAVAO OCCIDENTAL","REGION XI (DAVAO REGION)",39,120,189),
    @("DAVAO ORIENTAL","REGION XI (DAVAO REGION)",66,218,383),
    @("DINAGAT ISLANDS","REGION XIII (CARAGA)",234,36,139),
    @("EASTERN SAMAR","REGION VIII (EASTERN VISAYAS)",125,200,529),
    @("GUIMARAS","REGION VI (WESTERN VISAYAS)",47,95,114),
    @("IFUGAO","CORDILLERA ADMINISTRATIVE REGION (CAR)",197,212,263),
    @("ILOCOS NORTE","REGION I (ILOCOS REGION)",273,181,443),
    @("ILOCOS SUR","REGION I (ILOCOS REGION)",294,503,567),
    @("ILOILO","REGION VI (WESTERN VISAYAS)",386,660,1279),
    @("ISABELA","REGION II (CAGAYAN VALLEY)",468,964,1143),
    @("KALINGA","CORDILLERA ADMINISTRATIVE REGION (CAR)",179,158,292),
    @("LA UNION","REGION I (ILOCOS REGION)",216,346,441),
    @("LAGUNA","REGION IV-A (CALABARZON)",486,610,591),
    @("LANAO DEL NORTE","REGION X (NORTHERN MINDANAO)",137,303,489),
    @("LANAO DEL SUR","AUTONOMOUS REGION IN MUSLIM MINDANAO (ARMM)",0,242,879),
    @("LEYTE","REGION VIII (EASTERN VISAYAS)",11,605,1500),
    @("MAGUINDANAO","AUTONOMOUS REGION IN MUSLIM MINDANAO (ARMM)",0,344,633),
    @("MARINDUQUE","MIMAROPA REGION",119,48,228),
    @("MASBATE","REGION V (BICOL REGION)",202,383,737),
    @("MISAMIS OCCIDENTAL","REGION X (NORTHERN MINDANAO)",95,183,531),
    @("MISAMIS ORIENTAL","REGION X (NORTHERN MINDANAO)",171,519,651),
    @("MOUNTAIN PROVINCE","CORDILLERA ADMINISTRATIVE REGION (CAR)",85,149,264),
    @("NCR, CITY OF MANILA, FIRST DISTRICT (NOT A PROVINCE)","NATIONAL CAPITAL REGION (NCR)",0,0,108),
    @("NCR, FOURTH DISTRICT (NOT A PROVINCE)","NATIONAL CAPITAL REGION (NCR)",212,0,217),
    @("NCR, SECOND DISTRICT (NOT A PROVINCE)","NATIONAL CAPITAL REGION (NCR)",139,0,255),
    @("NCR, THIRD DISTRICT (NOT A PROVINCE)","NATIONAL CAPITAL REGION (NCR)",95,0,206),
    @("NEGROS OCCIDENTAL","REGION VI (WESTERN VISAYAS)",574,764,1207),
    @("NEGROS ORIENTAL","REGION VII (CENTRAL VISAYAS)",290,552,894),
    @("NORTHERN SAMAR","REGION VIII (EASTERN VISAYAS)",109,231,603),
    @("NUEVA ECIJA","REGION III (CENTRAL LUZON)",278,511,935),
    @("NUEVA VIZCAYA","REGION II (CAGAYAN VALLEY)",231,268,379),
    @("OCCIDENTAL MINDORO","MIMAROPA REGION",84,151,348),
    @("ORIENTAL MINDORO","MIMAROPA REGION",579,465,566),
    @("PALAWAN","MIMAROPA REGION",401,430,875),
    @("PAMPANGA","REGION III (CENTRAL LUZON)",286,650,727),
    @("PANGASINAN","REGION I (ILOCOS REGION)",387,1195,1505),
    @("QUEZON","REGION IV-A (CALABARZON)",911,869,1047),
    @("QUIRINO","REGION II (CAGAYAN VALLEY)",25,159,211),
    @("RIZAL","REGION IV-A (CALABARZON)",331,398,353),
    @("ROMBLON","MIMAROPA REGION",289,221,260),
    @("SAMAR (WESTERN SAMAR)","REGION VIII (EASTERN VISAYAS)",100,253,983),
    @("SARANGANI","REGION XII (SOCCSKSARGEN)",150,170,405),
    @("SIQUIJOR","REGION VII (CENTRAL VISAYAS)",94,67,79),
    @("SORSOGON","REGION V (BICOL REGION)",185,450,615),
    @("SOUTH COTABATO","REGION XII (SOCCSKSARGEN)",390,369,543),
    @("SOUTHERN LEYTE","REGION VIII (EASTERN VISAYAS)",525,175,418),
    @("SULTAN KUDARAT","REGION XII (SOCCSKSARGEN)",81,302,442),
    @("SULU","AUTONOMOUS REGION IN MUSLIM MINDANAO (ARMM)",0,171,456),
    @("SURIGAO DEL NORTE","REGION XIII (CARAGA)",224,201,428),
    @("SURIGAO DEL SUR","REGION XIII (CARAGA)",362,313,570),
    @("TARLAC","REGION III (CENTRAL LUZON)",194,309,599),
    @("TAWI-TAWI","AUTONOMOUS REGION IN MUSLIM MINDANAO (ARMM)",0,103,262),
    @("ZAMBALES","REGION III (CENTRAL LUZON)",124,274,368),
    @("ZAMBOANGA DEL NORTE","REGION IX (ZAMBOANGA PENINSULA)",178,431,865),
    @("ZAMBOANGA DEL SUR","REGION IX (ZAMBOANGA PENINSULA)",293,434,1097),
    @("ZAMBOANGA SIBUGAY","REGION IX (ZAMBOANGA PENINSULA)",124,170,506)
)

$startRow = 27
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

$ws.Range("A1:E88").Select()
